# Auto-generated Excel COM-interop script to update cryptos.xlsx
# Applies the cell-value changes described by the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper cell used as a scratch pad to inject literal text without
# Excel re-interpreting numeric-looking strings (e.g. "1.00", "599.03")
# as actual numbers. We write a string-literal formula into it, copy the
# *computed value* and paste-special (values only) into the target cell,
# which stores the exact text without touching any cell style.
$scratch = $ws.Range("Z1")

$scratch.Formula = "=""67.360.19"""
$scratch.Copy()
$ws.Range("D2").PasteSpecial(-4163)
$scratch.Formula = "=""  +1.05%  """
$scratch.Copy()
$ws.Range("E2").PasteSpecial(-4163)
$scratch.Formula = "=""3.509.94"""
$scratch.Copy()
$ws.Range("D3").PasteSpecial(-4163)
$scratch.Formula = "=""  +0.32%  """
$scratch.Copy()
$ws.Range("E3").PasteSpecial(-4163)
$scratch.Formula = "=""1.00"""
$scratch.Copy()
$ws.Range("D4").PasteSpecial(-4163)
$scratch.Formula = "=""  +0.02%  """
$scratch.Copy()
$ws.Range("E4").PasteSpecial(-4163)
$scratch.Formula = "=""599.03"""
$scratch.Copy()
$ws.Range("D5").PasteSpecial(-4163)
$scratch.Formula = "=""  +1.00%  """
$scratch.Copy()
$ws.Range("E5").PasteSpecial(-4163)
$scratch.Formula = "=""176.30"""
$scratch.Copy()
$ws.Range("D6").PasteSpecial(-4163)
$scratch.Formula = "=""  +4.31%  """
$scratch.Copy()
$ws.Range("E6").PasteSpecial(-4163)
$scratch.Formula = "=""1.00"""
$scratch.Copy()
$ws.Range("D7").PasteSpecial(-4163)
$scratch.Formula = "=""  +0.02%  """
$scratch.Copy()
$ws.Range("E7").PasteSpecial(-4163)
$scratch.Formula = "=""0.587"""
$scratch.Copy()
$ws.Range("D8").PasteSpecial(-4163)
$scratch.Formula = "=""  -0.51%  """
$scratch.Copy()
$ws.Range("E8").PasteSpecial(-4163)
$scratch.Formula = "=""  -0.32%  """
$scratch.Copy()
$ws.Range("E9").PasteSpecial(-4163)
$scratch.Formula = "=""  -1.90%  """
$scratch.Copy()
$ws.Range("E10").PasteSpecial(-4163)
$scratch.Formula = "=""  -0.01%  """
$scratch.Copy()
$ws.Range("E11").PasteSpecial(-4163)
$scratch.Formula = "=""4.117.24"""
$scratch.Copy()
$ws.Range("D12").PasteSpecial(-4163)
$scratch.Formula = "=""  +0.28%  """
$scratch.Copy()
$ws.Range("E12").PasteSpecial(-4163)
$scratch.Formula = "=""30.77"""
$scratch.Copy()
$ws.Range("D13").PasteSpecial(-4163)
$scratch.Formula = "=""  +8.73%  """
$scratch.Copy()
$ws.Range("E13").PasteSpecial(-4163)
$scratch.Formula = "=""  +0.21%  """
$scratch.Copy()
$ws.Range("E14").PasteSpecial(-4163)
$scratch.Formula = "=""67.376.54"""
$scratch.Copy()
$ws.Range("D15").PasteSpecial(-4163)
$scratch.Formula = "=""  +0.99%  """
$scratch.Copy()
$ws.Range("E15").PasteSpecial(-4163)
$scratch.Formula = "=""  -1.00%  """
$scratch.Copy()
$ws.Range("E16").PasteSpecial(-4163)
$scratch.Formula = "=""3.493.29"""
$scratch.Copy()
$ws.Range("D17").PasteSpecial(-4163)
$scratch.Formula = "=""  -0.17%  """
$scratch.Copy()
$ws.Range("E17").PasteSpecial(-4163)
$scratch.Formula = "=""  -0.05%  """
$scratch.Copy()
$ws.Range("E18").PasteSpecial(-4163)
$scratch.Formula = "=""14.59"""
$scratch.Copy()
$ws.Range("D19").PasteSpecial(-4163)
$scratch.Formula = "=""  +4.13%  """
$scratch.Copy()
$ws.Range("E19").PasteSpecial(-4163)
$scratch.Formula = "=""394.67"""
$scratch.Copy()
$ws.Range("D20").PasteSpecial(-4163)
$scratch.Formula = "=""  -0.44%  """
$scratch.Copy()
$ws.Range("E20").PasteSpecial(-4163)
$scratch.Formula = "=""  +0.56%  """
$scratch.Copy()
$ws.Range("E21").PasteSpecial(-4163)
$scratch.Formula = "=""73.58"""
$scratch.Copy()
$ws.Range("D22").PasteSpecial(-4163)
$scratch.Formula = "=""  +0.08%  """
$scratch.Copy()
$ws.Range("E22").PasteSpecial(-4163)
$scratch.Formula = "=""  +0.20%  """
$scratch.Copy()
$ws.Range("E23").PasteSpecial(-4163)
$scratch.Formula = "=""  +0.71%  """
$scratch.Copy()
$ws.Range("E24").PasteSpecial(-4163)
$scratch.Formula = "=""  -0.63%  """
$scratch.Copy()
$ws.Range("E25").PasteSpecial(-4163)
$scratch.Formula = "=""  +0.47%  """
$scratch.Copy()
$ws.Range("E26").PasteSpecial(-4163)
$scratch.Formula = "=""10.22"""
$scratch.Copy()
$ws.Range("D27").PasteSpecial(-4163)
$scratch.Formula = "=""  +0.63%  """
$scratch.Copy()
$ws.Range("E27").PasteSpecial(-4163)
$scratch.Formula = "=""  +0.08%  """
$scratch.Copy()
$ws.Range("E28").PasteSpecial(-4163)
$scratch.Formula = "=""  -0.54%  """
$scratch.Copy()
$ws.Range("E29").PasteSpecial(-4163)
$scratch.Formula = "=""6.18"""
$scratch.Copy()
$ws.Range("D30").PasteSpecial(-4163)
$scratch.Formula = "=""  -2.15%  """
$scratch.Copy()
$ws.Range("E30").PasteSpecial(-4163)
$scratch.Formula = "=""  -1.96%  """
$scratch.Copy()
$ws.Range("E31").PasteSpecial(-4163)
$scratch.Formula = "=""  -0.02%  """
$scratch.Copy()
$ws.Range("E32").PasteSpecial(-4163)
$scratch.Formula = "=""  -0.33%  """
$scratch.Copy()
$ws.Range("E33").PasteSpecial(-4163)
$scratch.Formula = "=""7.40"""
$scratch.Copy()
$ws.Range("D34").PasteSpecial(-4163)
$scratch.Formula = "=""  +0.12%  """
$scratch.Copy()
$ws.Range("E34").PasteSpecial(-4163)
$scratch.Formula = "=""1.65"""
$scratch.Copy()
$ws.Range("D35").PasteSpecial(-4163)
$scratch.Formula = "=""  +2.14%  """
$scratch.Copy()
$ws.Range("E35").PasteSpecial(-4163)
$scratch.Formula = "=""164.21"""
$scratch.Copy()
$ws.Range("D36").PasteSpecial(-4163)
$scratch.Formula = "=""  +1.08%  """
$scratch.Copy()
$ws.Range("E36").PasteSpecial(-4163)
$scratch.Formula = "=""0.881"""
$scratch.Copy()
$ws.Range("D37").PasteSpecial(-4163)
$scratch.Formula = "=""  -1.95%  """
$scratch.Copy()
$ws.Range("E37").PasteSpecial(-4163)
$scratch.Formula = "=""1.92"""
$scratch.Copy()
$ws.Range("D38").PasteSpecial(-4163)
$scratch.Formula = "=""  +0.77%  """
$scratch.Copy()
$ws.Range("E38").PasteSpecial(-4163)
$scratch.Formula = "=""7.05"""
$scratch.Copy()
$ws.Range("D39").PasteSpecial(-4163)
$scratch.Formula = "=""  +3.57%  """
$scratch.Copy()
$ws.Range("E39").PasteSpecial(-4163)
$scratch.Formula = "=""27.73"""
$scratch.Copy()
$ws.Range("D40").PasteSpecial(-4163)
$scratch.Formula = "=""  +2.13%  """
$scratch.Copy()
$ws.Range("E40").PasteSpecial(-4163)
$scratch.Formula = "=""4.69"""
$scratch.Copy()
$ws.Range("D41").PasteSpecial(-4163)
$scratch.Formula = "=""  -0.18%  """
$scratch.Copy()
$ws.Range("E41").PasteSpecial(-4163)
$scratch.Formula = "=""  -1.36%  """
$scratch.Copy()
$ws.Range("E42").PasteSpecial(-4163)
$scratch.Formula = "=""26.25"""
$scratch.Copy()
$ws.Range("D43").PasteSpecial(-4163)
$scratch.Formula = "=""  -0.78%  """
$scratch.Copy()
$ws.Range("E43").PasteSpecial(-4163)
$scratch.Formula = "=""2.805.74"""
$scratch.Copy()
$ws.Range("D44").PasteSpecial(-4163)
$scratch.Formula = "=""  +0.31%  """
$scratch.Copy()
$ws.Range("E44").PasteSpecial(-4163)
$scratch.Formula = "=""dogwifhat"""
$scratch.Copy()
$ws.Range("B45").PasteSpecial(-4163)
$scratch.Formula = "=""https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"""
$scratch.Copy()
$ws.Range("C45").PasteSpecial(-4163)
$scratch.Formula = "=""2.56"""
$scratch.Copy()
$ws.Range("D45").PasteSpecial(-4163)
$scratch.Formula = "=""  -0.14%  """
$scratch.Copy()
$ws.Range("E45").PasteSpecial(-4163)
$scratch.Formula = "=""OKB"""
$scratch.Copy()
$ws.Range("B46").PasteSpecial(-4163)
$scratch.Formula = "=""https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"""
$scratch.Copy()
$ws.Range("C46").PasteSpecial(-4163)
$scratch.Formula = "=""42.54"""
$scratch.Copy()
$ws.Range("D46").PasteSpecial(-4163)
$scratch.Formula = "=""  -0.68%  """
$scratch.Copy()
$ws.Range("E46").PasteSpecial(-4163)
$scratch.Formula = "=""  -2.86%  """
$scratch.Copy()
$ws.Range("E47").PasteSpecial(-4163)
$scratch.Formula = "=""343.29"""
$scratch.Copy()
$ws.Range("D48").PasteSpecial(-4163)
$scratch.Formula = "=""  +0.42%  """
$scratch.Copy()
$ws.Range("E48").PasteSpecial(-4163)
$scratch.Formula = "=""  -0.78%  """
$scratch.Copy()
$ws.Range("E49").PasteSpecial(-4163)
$scratch.Formula = "=""33.68"""
$scratch.Copy()
$ws.Range("D50").PasteSpecial(-4163)
$scratch.Formula = "=""  -0.79%  """
$scratch.Copy()
$ws.Range("E50").PasteSpecial(-4163)
$scratch.Formula = "=""  -1.11%  """
$scratch.Copy()
$ws.Range("E51").PasteSpecial(-4163)

$scratch.Clear()
$excel.CutCopyMode = 0

